$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.001461982727051
$ws.Range("B1").Value = 1.642426729202271
$ws.Range("C1").Value = 4.239952087402344
$ws.Range("D1").Value = 2.478299140930176
$ws.Range("E1").Value = 1.387456893920898
